$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel sending-cluster "Inflammatory-Mac" -> "MuSCs" (rows 8-10)
# (the workbook already has a "MuSCs" target-cluster label elsewhere;
#  after this rename the two converge and the old duplicate shared string is dropped)
$ws.Range("A8").Value = "MuSCs"
$ws.Range("A9").Value = "MuSCs"
$ws.Range("A10").Value = "MuSCs"

# Updated TPM-derived NATMI metrics (columns E:T) for every data row
# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2955753333333333
$ws.Range("H2").Value = 0.8867259999999999
$ws.Range("I2").Value = 0.239018529794766
$ws.Range("J2").Value = 0.2584571780171812
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1882183333333333
$ws.Range("N2").Value = 0.564655
$ws.Range("O2").Value = 0.01184593174728904
$ws.Range("P2").Value = 0.01206006386170697
$ws.Range("Q2").Value = 0.05563269661444444
$ws.Range("R2").Value = 0.50069426953
$ws.Range("S2").Value = 0.00283139719028617
$ws.Range("T2").Value = 0.003117010072403771

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2955753333333333
$ws.Range("H3").Value = 0.8867259999999999
$ws.Range("I3").Value = 0.239018529794766
$ws.Range("J3").Value = 0.2584571780171812
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 14.85429766666667
$ws.Range("N3").Value = 44.562893
$ws.Range("O3").Value = 0.9348876551872286
$ws.Range("P3").Value = 0.9517870831612478
$ws.Range("Q3").Value = 4.390563984257556
$ws.Range("R3").Value = 39.515075858318
$ws.Range("S3").Value = 0.2234554728661275
$ws.Range("T3").Value = 0.2459962035870602

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2955753333333333
$ws.Range("H4").Value = 0.8867259999999999
$ws.Range("I4").Value = 0.239018529794766
$ws.Range("J4").Value = 0.2584571780171812
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.8463425
$ws.Range("N4").Value = 1.692685
$ws.Range("O4").Value = 0.05326641306548233
$ws.Range("P4").Value = 0.0361528529770452
$ws.Range("Q4").Value = 0.2501579665516667
$ws.Range("R4").Value = 1.50094779931
$ws.Range("S4").Value = 0.0127316597383523
$ws.Range("T4").Value = 0.009343964357717149

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.258422
$ws.Range("H5").Value = 0.775266
$ws.Range("I5").Value = 0.2089742936599006
$ws.Range("J5").Value = 0.2259695357671569
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1882183333333333
$ws.Range("N5").Value = 0.564655
$ws.Range("O5").Value = 0.01184593174728904
$ws.Range("P5").Value = 0.01206006386170697
$ws.Range("Q5").Value = 0.04863975813666667
$ws.Range("R5").Value = 0.43775782323
$ws.Range("S5").Value = 0.00247549521963312
$ws.Range("T5").Value = 0.00272520703215219

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.258422
$ws.Range("H6").Value = 0.775266
$ws.Range("I6").Value = 0.2089742936599006
$ws.Range("J6").Value = 0.2259695357671569
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 14.85429766666667
$ws.Range("N6").Value = 44.562893
$ws.Range("O6").Value = 0.9348876551872286
$ws.Range("P6").Value = 0.9517870831612478
$ws.Range("Q6").Value = 3.838677311615333
$ws.Range("R6").Value = 34.548095804538
$ws.Range("S6").Value = 0.1953674873941118
$ws.Range("T6").Value = 0.2150748853311236

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.258422
$ws.Range("H7").Value = 0.775266
$ws.Range("I7").Value = 0.2089742936599006
$ws.Range("J7").Value = 0.2259695357671569
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.8463425
$ws.Range("N7").Value = 1.692685
$ws.Range("O7").Value = 0.05326641306548233
$ws.Range("P7").Value = 0.0361528529770452
$ws.Range("Q7").Value = 0.218713521535
$ws.Range("R7").Value = 1.31228112921
$ws.Range("S7").Value = 0.01113131104615567
$ws.Range("T7").Value = 0.008169443403881182

# Row 8
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 0.27902
$ws.Range("H8").Value = 0.55804
$ws.Range("I8").Value = 0.2256309734348681
$ws.Range("J8").Value = 0.1626539016795581
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.1882183333333333
$ws.Range("N8").Value = 0.564655
$ws.Range("O8").Value = 0.01184593174728904
$ws.Range("P8").Value = 0.01206006386170697
$ws.Range("Q8").Value = 0.05251667936666667
$ws.Range("R8").Value = 0.3151000762
$ws.Range("S8").Value = 0.002672809111383834
$ws.Range("T8").Value = 0.001961616441611276

# Row 9
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.5
$ws.Range("G9").Value = 0.27902
$ws.Range("H9").Value = 0.55804
$ws.Range("I9").Value = 0.2256309734348681
$ws.Range("J9").Value = 0.1626539016795581
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 14.85429766666667
$ws.Range("N9").Value = 44.562893
$ws.Range("O9").Value = 0.9348876551872286
$ws.Range("P9").Value = 0.9517870831612478
$ws.Range("Q9").Value = 4.144646134953334
$ws.Range("R9").Value = 24.86787680972
$ws.Range("S9").Value = 0.2109396116921357
$ws.Range("T9").Value = 0.1548118826443829

# Row 10
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.5
$ws.Range("G10").Value = 0.27902
$ws.Range("H10").Value = 0.55804
$ws.Range("I10").Value = 0.2256309734348681
$ws.Range("J10").Value = 0.1626539016795581
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.8463425
$ws.Range("N10").Value = 1.692685
$ws.Range("O10").Value = 0.05326641306548233
$ws.Range("P10").Value = 0.0361528529770452
$ws.Range("Q10").Value = 0.23614648435
$ws.Range("R10").Value = 0.9445859374
$ws.Range("S10").Value = 0.01201855263134855
$ws.Range("T10").Value = 0.005880402593563828

# Row 11
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.4036036666666667
$ws.Range("H11").Value = 1.210811
$ws.Range("I11").Value = 0.3263762031104653
$ws.Range("J11").Value = 0.3529193845361038
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.1882183333333333
$ws.Range("N11").Value = 0.564655
$ws.Range("O11").Value = 0.01184593174728904
$ws.Range("P11").Value = 0.01206006386170697
$ws.Range("Q11").Value = 0.07596560946722224
$ws.Range("R11").Value = 0.6836904852050001
$ws.Range("S11").Value = 0.003866230225985917
$ws.Range("T11").Value = 0.00425623031553973

# Row 12
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.4036036666666667
$ws.Range("H12").Value = 1.210811
$ws.Range("I12").Value = 0.3263762031104653
$ws.Range("J12").Value = 0.3529193845361038
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 14.85429766666667
$ws.Range("N12").Value = 44.562893
$ws.Range("O12").Value = 0.9348876551872286
$ws.Range("P12").Value = 0.9517870831612478
$ws.Range("Q12").Value = 5.995249004024778
$ws.Range("R12").Value = 53.95724103622301
$ws.Range("S12").Value = 0.3051250832348535
$ws.Range("T12").Value = 0.3359041115986811

# Row 13
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.4036036666666667
$ws.Range("H13").Value = 1.210811
$ws.Range("I13").Value = 0.3263762031104653
$ws.Range("J13").Value = 0.3529193845361038
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.8463425
$ws.Range("N13").Value = 1.692685
$ws.Range("O13").Value = 0.05326641306548233
$ws.Range("P13").Value = 0.0361528529770452
$ws.Range("Q13").Value = 0.3415869362558334
$ws.Range("R13").Value = 2.049521617535
$ws.Range("S13").Value = 0.0173848896496258
$ws.Range("T13").Value = 0.01275904262188304
